$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handback identifiers / timestamps (this run's regenerated report).
# ---------------------------------------------------------------------------
$oldMd1 = "046f5f43-136a-40ec-b0a8-eb093b422b4b.md"
$newMd1 = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.md"

$oldMd2 = "fe7cc5c2-cea7-482f-8494-90acb7a68223.md"
$newMd2 = "ffffb99fb7df-7bc5-4f5a-a233-b7bb0c6e0563.md"

$newXlfZh = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.zh-cn.xlf"
$newXlfDe = "f64fd7e2-04ab-449d-ada8-fdaffe3197cb.c4ef44521985cc8052aa5530f95c3ba80f4971b4.de-de.xlf"

$zhHandoffTime  = "2016-03-13 21:13:27"
$zhHandbackTime = "2016-03-13 21:13:51"
$deHandoffTime  = "2016-03-13 21:13:33"
$deHandbackTime = "2016-03-13 21:13:57"

# ---------------------------------------------------------------------------
# Overview sheet: File Name column hyperlinks / cell text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd1) {
        $hl.TextToDisplay = $newMd1
    } elseif ($hl.TextToDisplay -eq $oldMd2) {
        $hl.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $zhHandoffTime
$wsZh.Range("F2").Value = $newMd1
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $zhHandbackTime

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $zhHandoffTime
$wsZh.Range("F3").Value = $newMd2
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $zhHandbackTime

$zhDisplays = @($newMd1, ".md", $newXlfZh, $newMd1, $newXlfZh, $newMd2, ".md", $newXlfZh, $newMd2, $newXlfZh)
$i = 0
foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = $zhDisplays[$i]
    $i = $i + 1
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $deHandoffTime
$wsDe.Range("F2").Value = $newMd1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $deHandbackTime

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $deHandoffTime
$wsDe.Range("F3").Value = $newMd2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $deHandbackTime

$deDisplays = @($newMd1, ".md", $newXlfDe, $newMd1, $newXlfDe, $newMd2, ".md", $newXlfDe, $newMd2, $newXlfDe)
$i = 0
foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = $deDisplays[$i]
    $i = $i + 1
}
